# Add battery status setpoints
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels in row 1
$ws.Range("AR1").Value = "CellTempMin"
$ws.Range("AS1").Value = "CellTempMax"
$ws.Range("AT1").Value = "CellTempRange"
$ws.Range("AU1").Value = "CellVoltageMin"
$ws.Range("AV1").Value = "CellVoltageMax"
$ws.Range("AW1").Value = "CellVoltageRange"

# New setpoint values in row 2. These look numeric, but the source data
# keeps them as literal text (e.g. "5.00"), so force a Text format before
# typing them in and then drop the format back so no visible formatting
# change is introduced.
$ws.Range("AR2:AW2").NumberFormat = "@"
$ws.Range("AR2").Value = "5.00"
$ws.Range("AS2").Value = "40.00"
$ws.Range("AT2").Value = "15.00"
$ws.Range("AU2").Value = "2.50"
$ws.Range("AV2").Value = "4.20"
$ws.Range("AW2").Value = "50.00"
$ws.Range("AR2:AW2").Style = "Normal"
